# Bugfixed QoQ Visualizations and a typo in the evaluation objects
#
# The source data table (date_of_forecast | y_0 | y_0_forecast | y_1 | y_1_forecast)
# contained four extra leading years (1984-1987) that should not have been part of
# the series. Remove those four data rows (rows 2-5) so the table starts with the
# 1988 forecast, shifting every remaining row up by four positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the four obsolete rows (old years 1984-1987); Excel automatically shifts
# the rows below upward, which reproduces the row-for-row remapping seen in the diff.
$ws.Range("A2:E5").EntireRow.Delete()
